$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptos list, refreshed by the scheduled GitHub Actions scraper run.
# Values are written as literal text (NumberFormat "@") because the sheet
# stores these figures as plain strings (e.g. "29.825.93", "1.000") rather
# than numbers, and some of them would otherwise be silently normalised
# (trailing zeros / thousands-style dots dropped) by Excel's automatic
# type coercion on assignment.

$rows = @(
    @{ Row = 2;  D = "29.825.93";   E = "  -0.14%  " }
    @{ Row = 3;  D = "1.877.84";    E = "  +0.34%  " }
    @{ Row = 4;  D = "1.002";       E = "  +0.23%  " }
    @{ Row = 5;  D = "0.7288";      E = "  -0.99%  " }
    @{ Row = 6;  D = "241.91";      E = "  -0.08%  " }
    @{ Row = 7;  D = "1.002";       E = "  +0.22%  " }
    @{ Row = 8;  D = "0.3135";      E = "  -0.97%  " }
    @{ Row = 9;  D = "0.07099";     E = "  +0.04%  " }
    @{ Row = 10; D = "24.39";       E = "  -1.10%  " }
    @{ Row = 11; D = "0.08277";     E = "  -1.86%  " }
    @{ Row = 12; D = "0.7476";      E = "  -0.42%  " }
    @{ Row = 13; D = "1.892.36";    E = "  +1.38%  " }
    @{ Row = 14; D = "5.329";       E = "  -0.80%  " }
    @{ Row = 15; D = "92.58";       E = "  +0.20%  " }
    @{ Row = 16; D = "29.859.84";   E = "  -0.01%  " }
    @{ Row = 17; D = "6.058";       E = "  +0.14%  " }
    @{ Row = 18; D = "248.22";      E = "  +2.09%  " }
    @{ Row = 19; D = "13.38";       E = "  -1.31%  " }
    @{ Row = 20; D = "0.000007834"; E = "  +0.19%  " }
    @{ Row = 21; D = "2.173.26";    E = "  +2.43%  " }
    @{ Row = 22; D = $null;         E = "  +0.29%  " }
    @{ Row = 23; D = "1.002";       E = "  +0.22%  " }
    @{ Row = 24; D = "7.745";       E = "  -2.48%  " }
    @{ Row = 25; D = "0.1541";      E = "  -1.79%  " }
    @{ Row = 26; D = "9.179";       E = "  -1.50%  " }
    @{ Row = 27; D = "163.13";      E = "  -0.56%  " }
    @{ Row = 28; D = "18.56";       E = "  -0.37%  " }
    @{ Row = 29; D = "2.030";       E = "  +0.14%  " }
    @{ Row = 30; D = "1.441";       E = "  -1.74%  " }
    @{ Row = 31; D = "4.548";       E = "  +0.16%  " }
    @{ Row = 32; D = "1.529";       E = "  -0.30%  " }
    @{ Row = 33; D = "4.193";       E = "  -1.48%  " }
    @{ Row = 34; D = "0.05265";     E = "  -1.23%  " }
    @{ Row = 35; D = "1.234";       E = "  +0.06%  " }
    @{ Row = 36; D = "0.7576";      E = "  +1.11%  " }
    @{ Row = 37; D = "1.000";       E = "  +0.12%  " }
    @{ Row = 38; D = "2.716";       E = "  +0.99%  " }
    @{ Row = 39; D = "0.01932";     E = "  -0.90%  " }
    @{ Row = 40; D = "2.754";       E = "  -0.12%  " }
    @{ Row = 41; D = "0.4494";      E = "  +0.58%  " }
    @{ Row = 42; D = "6.008";       E = "  -0.77%  " }
    @{ Row = 43; D = "0.8685";      E = "  +0.08%  " }
    @{ Row = 44; D = "71.40";       E = "  -1.13%  " }
    @{ Row = 45; D = "1.069.52";    E = "  -2.41%  " }
    @{ Row = 46; D = $null;         E = "  +2.23%  " }
    @{ Row = 47; D = $null;         E = "  +0.17%  " }
    @{ Row = 48; D = "1.834";       E = "  -0.19%  " }
    @{ Row = 49; D = "7.519";       E = "  -2.75%  " }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($null -ne $r.D) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
    }
    $cell = $ws.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = $r.E
}

# Rank #48 (row 50) and #49 (row 51) in the list changed this run: a new
# coin, EnergySwap, entered the ranking at #48 - pushing RocketPoolETH
# down to #49 with a refreshed price - while the previous #49
# (SynthetixNetwork) dropped off the bottom of the sheet.
$ws.Range("B50:E51").NumberFormat = "@"

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.515"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.048.18"
$ws.Range("E51").Value = "  +1.16%  "
